# Excel Group import: auto-enable Web Hook Results reporting for groups
# whose Title matches a 32-character hex string (CVDLS-220).
# Add a new sample "Individual Group" row demonstrating a matching
# (enabled-for-reporting) Group Title value in both the ServiceNow ID
# and Title columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = "09876543210987654321abcdefABCDEF"
$ws.Range("H7").Value = 1
$ws.Range("J7").Value = "09876543210987654321abcdefABCDEF"

$ws.Range("C7").Select()
